$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 152, shifting existing rows 152:356 down to 153:357.
$ws.Rows.Item(152).Insert()

# Populate the newly inserted row 152 with its data.
$ws.Range("A152").Value = 3
$ws.Range("B152").Value = "Femacal de La Calera"
$ws.Range("C152").Value = "Coquimbo"
$ws.Range("D152").Value = 44799
$ws.Range("D152").NumberFormat = $ws.Range("D153").NumberFormat
$ws.Range("E152").Value = 5
$ws.Range("F152").Value = 100112039
$ws.Range("G152").Value = "Ciboulette"
$ws.Range("H152").Value = "Sin especificar"
$ws.Range("I152").Value = "Primera"
$ws.Range("J152").Value = 160
$ws.Range("K152").Value = 1500
$ws.Range("L152").Value = 1500
$ws.Range("M152").Value = 1500
$ws.Range("N152").Value = "$/docena de atados"
$ws.Range("O152").Value = "Provincia de Quillota"
$ws.Range("P152").Value = 500
$ws.Range("Q152").Value = 3
$ws.Range("R152").Value = "Hortaliza"
